$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (row, A..O)
$rows = @(
    @(9,  '2023-12-15 01:16:26', 7,  6,  3, 2, 0, 1, 0.001, 0.01, 0.003, 100, 512, 10, 7, 0.8571428571428571),
    @(10, '2023-12-15 01:21:13', 3,  3,  1, 1, 0, 0, 0.001, 0.01, 0.003, 100, 512, 10, 7, 1),
    @(11, '2023-12-16 15:20:34', 28, 20, 5, 3, 8, 7, 0.001, 0.01, 0.003, 100, 512, 10, 7, 0.7142857142857143)
)

foreach ($row in $rows) {
    $r = $row[0]

    # Column A is stored as plain text, not a date value
    $ws.Cells.Item($r, 1).Value = [string]$row[1]

    for ($col = 2; $col -le 15; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col]
    }
}
